# Weekly update: prepend a new week of "Pepino dulce" price records.
# Insert 4 new rows right before row 351 (shifts existing rows 351-395
# down to 355-399, content unchanged) and populate the new rows with the
# new week's data (fecha = 45077).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert four whole rows above the current row 351.
$ws.Range("A351:A354").EntireRow.Insert()

# Common / constant columns shared by every record in this block.
$mercadoId   = 6
$mercado     = "Mercado Mayorista Lo Valledor de Santiago"
$region      = "Metropolitana"
$fecha       = 45077
$codreg      = 13
$categoriaId = 100112043
$categoria   = "Pepino dulce"
$variedad    = "Cultivar IV Región"
$unidad      = "`$/bandeja 18 kilos"
$origen      = "Provincia de Limarí"
$kgUnidades  = 18
$clasif      = "Hortaliza"

$rows = @(
    @{ Row = 351; Calidad = "Especial"; Volumen = 280; PMin = 14000; PMax = 14000; PProm = 14000; PKg = 778 },
    @{ Row = 352; Calidad = "Primera";  Volumen = 470; PMin = 12000; PMax = 12000; PProm = 12000; PKg = 667 },
    @{ Row = 353; Calidad = "Segunda";  Volumen = 300; PMin = 10000; PMax = 10000; PProm = 10000; PKg = 556 },
    @{ Row = 354; Calidad = "Tercera";  Volumen = 120; PMin = 8000;  PMax = 8000;  PProm = 8000;  PKg = 444 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $categoriaId
    $ws.Cells.Item($row, 7).Value  = $categoria
    $ws.Cells.Item($row, 8).Value  = $variedad
    $ws.Cells.Item($row, 9).Value  = $r.Calidad
    $ws.Cells.Item($row, 10).Value = $r.Volumen
    $ws.Cells.Item($row, 11).Value = $r.PMin
    $ws.Cells.Item($row, 12).Value = $r.PMax
    $ws.Cells.Item($row, 13).Value = $r.PProm
    $ws.Cells.Item($row, 14).Value = $unidad
    $ws.Cells.Item($row, 15).Value = $origen
    $ws.Cells.Item($row, 16).Value = $r.PKg
    $ws.Cells.Item($row, 17).Value = $kgUnidades
    $ws.Cells.Item($row, 18).Value = $clasif
}
